$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.3229990005493164
$ws.Range("A3").Value = 0.12165498733520508
$ws.Range("A4").Value = 0.11770105361938477
$ws.Range("A5").Value = 0.40745019912719727
$ws.Range("A6").Value = 1.6915290355682373
$ws.Range("A7").Value = 0.2897801399230957
$ws.Range("A8").Value = 0.17392897605895996
$ws.Range("A9").Value = 0.2842140197753906
$ws.Range("A10").Value = 0.43554091453552246
$ws.Range("A11").Value = 0.35520219802856445
$ws.Range("A12").Value = 1.7538189888000488
$ws.Range("A13").Value = 1.6837728023529053
$ws.Range("A14").Value = 1.4518780708312988
$ws.Range("A15").Value = 8.470750093460083
$ws.Range("A16").Value = 2.408931016921997
$ws.Range("A17").Value = 284.300626039505
$ws.Range("A18").Value = 410.2300720214844
$ws.Range("A19").Value = 3600.25691986084
$ws.Range("B19").Value = 0.000123033500004401
$ws.Range("C19").Value = 9326.000000000024
$ws.Range("E19").Value = 0.000123033500004401
$ws.Range("K19").Value = 9326.000000000024
$ws.Range("A20").Value = 75.52114486694336
$ws.Range("A21").Value = 1352.9148960113525
$ws.Range("B21").Value = 0.0000998685375930784
$ws.Range("C21").Value = 17145.00000000003
$ws.Range("D21").Value = "OPTIMAL"
$ws.Range("E21").Value = 0.0000998685375930784
$ws.Range("K21").Value = 17145.00000000003
$ws.Range("A22").Value = 3602.0627839565277
$ws.Range("A23").Value = 25.870085954666138
$ws.Range("A24").Value = 3600.0179359912872
$ws.Range("B24").Value = 0.0005101364049129629
$ws.Range("C24").Value = 6348.000000000327
$ws.Range("E24").Value = 0.0005101364049129629
$ws.Range("K24").Value = 6348.000000000327
$ws.Range("A25").Value = 3608.315566778183
$ws.Range("B25").Value = 0.006080916498046459
$ws.Range("C25").Value = 6207.0
$ws.Range("E25").Value = 0.006080916498046459
$ws.Range("K25").Value = 6207.0
$ws.Range("A26").Value = 1213.0529799461365
$ws.Range("A27").Value = 3601.960891008377
$ws.Range("B27").Value = 0.0008286443327412408
$ws.Range("C27").Value = 12435.0
$ws.Range("E27").Value = 0.0008286443327412408
$ws.Range("K27").Value = 12435.0
$ws.Range("A28").Value = 3604.9224960803986
$ws.Range("B28").Value = 0.002743763479497844
$ws.Range("E28").Value = 0.002743763479497844
$ws.Range("A29").Value = 3603.8255751132965
$ws.Range("B29").Value = 0.0004686373770336547
$ws.Range("C29").Value = 24970.0
$ws.Range("E29").Value = 0.0004686373770336547
$ws.Range("K29").Value = 24970.0
$ws.Range("A30").Value = 3603.810455083847
$ws.Range("B30").Value = 0.0012089126022245004
$ws.Range("C30").Value = 24588.0
$ws.Range("E30").Value = 0.0012089126022245004
$ws.Range("K30").Value = 24588.0
$ws.Range("A31").Value = 3605.4233560562134
$ws.Range("B31").Value = 0.004375204448806894
$ws.Range("C31").Value = 24456.0
$ws.Range("E31").Value = 0.004375204448806894
$ws.Range("K31").Value = 24456.0
$ws.Range("A32").Value = 3605.3445649147034
$ws.Range("B32").Value = 0.00033017274456214685
$ws.Range("C32").Value = 55420.00000000957
$ws.Range("E32").Value = 0.00033017274456214685
$ws.Range("K32").Value = 55420.00000000957
$ws.Range("A33").Value = 3602.5572600364685
$ws.Range("B33").Value = 0.0008888279143038277
$ws.Range("C33").Value = 54880.0
$ws.Range("E33").Value = 0.0008888279143038277
$ws.Range("K33").Value = 54880.0
$ws.Range("A34").Value = 3607.351403951645
$ws.Range("B34").Value = 0.003379859690148674
$ws.Range("C34").Value = 54736.0
$ws.Range("E34").Value = 0.003379859690148674
$ws.Range("K34").Value = 54736.0
$ws.Range("A35").Value = 3604.7848250865936
$ws.Range("B35").Value = 0.0003357184650025423
$ws.Range("C35").Value = 97855.0
$ws.Range("E35").Value = 0.0003357184650025423
$ws.Range("K35").Value = 97855.0
$ws.Range("A36").Value = 3607.671318054199
$ws.Range("B36").Value = 0.0014499311025643292
$ws.Range("C36").Value = 97246.0
$ws.Range("E36").Value = 0.0014499311025643292
$ws.Range("K36").Value = 97246.0
$ws.Range("A37").Value = 3607.0605511665344
$ws.Range("B37").Value = 0.002138993439047472
$ws.Range("C37").Value = 97242.0
$ws.Range("E37").Value = 0.002138993439047472
$ws.Range("K37").Value = 97242.0
$ws.Range("A38").Value = 2.037170886993408
$ws.Range("A39").Value = 8.591979026794434
$ws.Range("A40").Value = 8.96023416519165
$ws.Range("A41").Value = 1.3177781105041504
$ws.Range("A42").Value = 11.044636964797974
$ws.Range("A43").Value = 6.7002809047698975
$ws.Range("A44").Value = 10.944949865341187
$ws.Range("A45").Value = 9.37498688697815
$ws.Range("A46").Value = 184.60556602478027
$ws.Range("A47").Value = 3.841671943664551
$ws.Range("A48").Value = 18.37431001663208
$ws.Range("A49").Value = 333.44506096839905
$ws.Range("A50").Value = 24.974233150482178
$ws.Range("A51").Value = 40.21325898170471
$ws.Range("A52").Value = 250.65193891525269

Write-Host "applied changes"
